# LH_WF_USERHOME_REVIEWS.xlsx - "v1.1 close reviewer varification"
#
# Summary of the change being applied:
#  1. On the "LH-REVIEW-WF-USERHOME-SHEET" sheet, the reviewer verification
#     status in I2 moves from "open" to "closed".
#  2. On the "VERSION-HISTORY" sheet, a new version-history entry is added
#     in row 3 (v1.1 / Hala Eldaly / "close reviewer varification" /
#     29-Apr-2025), formatted the same way as the existing row 2 entry.
#  3. The remembered cell selections on both sheets are updated to match
#     where the author had clicked last (I19 on sheet 1, J9 on sheet 2),
#     leaving sheet 1 as the active tab.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet 2: VERSION-HISTORY
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VERSION-HISTORY")

# Duplicate the formatting of the existing v1.0 row onto the new row so the
# new entry looks like the rest of the table, then fill in its values.
$ws2.Range("A2:D2").Copy()
$ws2.Range("A3:D3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws2.Range("A3").Value = "v1.1"
$ws2.Range("B3").Value = "Hala Eldaly"
$ws2.Range("C3").Value = "close reviewer varification"
$ws2.Range("D3").Value = 45776

# ------------------------------------------------------------------
# Sheet 1: LH-REVIEW-WF-USERHOME-SHEET
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LH-REVIEW-WF-USERHOME-SHEET")

# "Reviewer verification" (column I) for the first review row flips from
# "open" to "closed".
$ws1.Range("I2").Value = "closed"

# ------------------------------------------------------------------
# Restore the saved cell selections / active sheet.
# ------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("I19").Select()

$ws2.Activate()
$ws2.Range("J9").Select()

$ws1.Activate()
